$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 2).Value = 'Bitcoin'
$ws.Cells.Item(2, 3).Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell 2 4 '29.988.56'
$ws.Cells.Item(2, 5).Value = '  +1.44%  '
$ws.Cells.Item(3, 2).Value = 'Ethereum'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell 3 4 '1.866.05'
$ws.Cells.Item(3, 5).Value = '  +1.21%  '
$ws.Cells.Item(4, 2).Value = 'TetherUSD'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 4 4 '0.9994'
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 5 4 '245.00'
$ws.Cells.Item(5, 5).Value = '  +0.42%  '
$ws.Cells.Item(6, 2).Value = 'XRP'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 6 4 '0.6543'
$ws.Cells.Item(6, 5).Value = '  +3.96%  '
$ws.Cells.Item(7, 2).Value = 'USDC'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 7 4 '1.0000'
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 2).Value = 'OKB'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 8 4 '47.92'
$ws.Cells.Item(8, 5).Value = '  +4.02%  '
$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 9 4 '0.07601'
$ws.Cells.Item(9, 5).Value = '  +2.11%  '
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 10 4 '0.2979'
$ws.Cells.Item(10, 5).Value = '  +0.52%  '
$ws.Cells.Item(11, 2).Value = 'Solana'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 11 4 '24.77'
$ws.Cells.Item(11, 5).Value = '  +5.45%  '
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 12 4 '0.07673'
$ws.Cells.Item(12, 5).Value = '  -0.10%  '
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 13 4 '1.875.71'
$ws.Cells.Item(13, 5).Value = '  +1.67%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 14 4 '5.094'
$ws.Cells.Item(14, 5).Value = '  +1.16%  '
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 15 4 '0.6977'
$ws.Cells.Item(15, 5).Value = '  +2.57%  '
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 16 4 '84.25'
$ws.Cells.Item(16, 5).Value = '  +0.76%  '
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 17 4 '0.000009754'
$ws.Cells.Item(17, 5).Value = '  +6.58%  '
$ws.Cells.Item(18, 2).Value = 'Uniswap'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 18 4 '6.176'
$ws.Cells.Item(18, 5).Value = '  +3.95%  '
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 19 4 '29.992.30'
$ws.Cells.Item(19, 5).Value = '  +1.53%  '
$ws.Cells.Item(20, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 20 4 '2.121.41'
$ws.Cells.Item(20, 5).Value = '  +1.03%  '
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 21 4 '238.08'
$ws.Cells.Item(21, 5).Value = '  -3.77%  '
$ws.Cells.Item(22, 2).Value = 'Avalanche'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 22 4 '12.76'
$ws.Cells.Item(22, 5).Value = '  +1.30%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 23 4 '0.9993'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '
$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 24 4 '7.782'
$ws.Cells.Item(24, 5).Value = '  +4.74%  '
$ws.Cells.Item(25, 2).Value = 'BinanceUSD'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 25 4 '1.001'
$ws.Cells.Item(25, 5).Value = '  -0.03%  '
$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 26 4 '0.1447'
$ws.Cells.Item(26, 5).Value = '  +2.14%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 27 4 '159.25'
$ws.Cells.Item(27, 5).Value = '  +0.42%  '
$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 28 4 '8.637'
$ws.Cells.Item(28, 5).Value = '  +0.56%  '
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 29 4 '18.02'
$ws.Cells.Item(29, 5).Value = '  +0.87%  '
$ws.Cells.Item(30, 2).Value = 'Hedera'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 30 4 '0.06102'
$ws.Cells.Item(30, 5).Value = '  +1.02%  '
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 31 4 '1.499'
$ws.Cells.Item(31, 5).Value = '  -0.06%  '
$ws.Cells.Item(32, 2).Value = 'Toncoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 32 4 '1.291'
$ws.Cells.Item(32, 5).Value = '  +5.23%  '
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 33 4 '4.186'
$ws.Cells.Item(33, 5).Value = '  +1.13%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 34 4 '4.126'
$ws.Cells.Item(34, 5).Value = '  +0.29%  '
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 35 4 '1.892'
$ws.Cells.Item(35, 5).Value = '  +0.53%  '
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 36 4 '1.180'
$ws.Cells.Item(36, 5).Value = '  +2.77%  '
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 37 4 '0.7396'
$ws.Cells.Item(37, 5).Value = '  +1.08%  '
$ws.Cells.Item(38, 2).Value = 'HuobiToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 38 4 '2.608'
$ws.Cells.Item(38, 5).Value = '  -0.17%  '
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 39 4 '2.825'
$ws.Cells.Item(39, 5).Value = '  -2.38%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 40 4 '0.01803'
$ws.Cells.Item(40, 5).Value = '  +1.87%  '
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 41 4 '1.219.47'
$ws.Cells.Item(41, 5).Value = '  -1.07%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 42 4 '6.399'
$ws.Cells.Item(42, 5).Value = '  +1.07%  '
$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 43 4 '0.9202'
$ws.Cells.Item(43, 5).Value = '  +0.21%  '
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 44 4 '7.916'
$ws.Cells.Item(44, 5).Value = '  +18.13%  '
$ws.Cells.Item(45, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 45 4 '2.031.85'
$ws.Cells.Item(45, 5).Value = '  +0.87%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 46 4 '1.000'
$ws.Cells.Item(46, 5).Value = '  -0.06%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 47 4 '67.69'
$ws.Cells.Item(47, 5).Value = '  +2.59%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 48 4 '101.96'
$ws.Cells.Item(48, 5).Value = '  -0.23%  '
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 49 4 '0.00000000123'
$ws.Cells.Item(49, 5).Value = '  +2.74%  '
$ws.Cells.Item(50, 2).Value = 'TheSandbox'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 50 4 '0.4106'
$ws.Cells.Item(50, 5).Value = '  +0.81%  '
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 51 4 '9.193'
$ws.Cells.Item(51, 5).Value = '  -0.99%  '